$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the execution time for "c++" (row 3, column C) from "1h 2m" to "1h 27m"
$ws.Range("C3").Value = "1h 27m"
